$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.560.50"
$ws.Range("E2").Value = "  +3.96%  "
$ws.Range("D3").Value = "1.846.09"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.031"
$ws.Range("E4").Value = "  +2.89%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.44"
$ws.Range("E5").Value = "  +4.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.026"
$ws.Range("E6").Value = "  +2.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4369"
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3743"
$ws.Range("E8").Value = "  +3.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07377"
$ws.Range("E9").Value = "  +2.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8755"
$ws.Range("E10").Value = "  +2.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.46"
$ws.Range("E11").Value = "  +4.11%  "
$ws.Range("D12").Value = "1.856.12"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.484"
$ws.Range("E13").Value = "  +3.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.680"
$ws.Range("E14").Value = "  +2.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07143"
$ws.Range("E15").Value = "  +3.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.67"
$ws.Range("E16").Value = "  +3.61%  "
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009004"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("E19").Value = "  +2.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.39"
$ws.Range("E20").Value = "  +2.54%  "
$ws.Range("D21").Value = "27.564.43"
$ws.Range("E21").Value = "  +3.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.258"
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.19"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").Value = "2.064.05"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.43"
$ws.Range("E25").Value = "  +3.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.930"
$ws.Range("E26").Value = "  +5.76%  "
$ws.Range("E27").Value = "  +3.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.249"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.943"
$ws.Range("E29").Value = "  +2.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.01"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09074"
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.203"
$ws.Range("E32").Value = "  +5.29%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7670"
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.497"
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("E35").Value = "  +4.81%  "
$ws.Range("E36").Value = "  +2.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.143"
$ws.Range("E37").Value = "  +2.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01974"
$ws.Range("E38").Value = "  +3.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05257"
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5165"
$ws.Range("E40").Value = "  +3.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.804"
$ws.Range("E41").Value = "  +7.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1672"
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.661"
$ws.Range("E43").Value = "  +3.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.533"
$ws.Range("E44").Value = "  +3.53%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.63"
$ws.Range("E45").Value = "  +3.44%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "108.69"
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.714"
$ws.Range("E47").Value = "  +4.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4648"
$ws.Range("E48").Value = "  +2.73%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06367"
$ws.Range("E49").Value = "  +2.59%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.895"
$ws.Range("E50").Value = "  +7.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.40"
$ws.Range("E51").Value = "  +6.43%  "
